$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated FilesTab query (B4): the participant/sample id fallbacks were
# simplified from 'Not specified in data' to '' when the CDS test suites
# were merged into the CDS_Regression suite.
$newFilesTabQuery = @'
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE COALESCE(g.library_strategy, "Not specified in data") in ['Not specified in data']
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name,'') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id, '') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER BY f.file_name limit 100
'@

$ws.Range("B4").Value = $newFilesTabQuery

# Update the selection to match the saved workbook state (B2 selected).
$ws.Range("B2").Select()
